$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 16) of data, mirroring the structure of the existing rows.
$row = 16

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42619.890405092592

$ws.Cells.Item($row, 2).Value = 4
$ws.Cells.Item($row, 3).Value = 57
$ws.Cells.Item($row, 4).Value = 40
$ws.Cells.Item($row, 5).Value = 57
$ws.Cells.Item($row, 6).Value = 47
$ws.Cells.Item($row, 7).Value = 22512
$ws.Cells.Item($row, 8).Value = 22157
$ws.Cells.Item($row, 9).Value = 1223
$ws.Cells.Item($row, 10).Value = 273
$ws.Cells.Item($row, 11).Value = 193
$ws.Cells.Item($row, 12).Value = 22
$ws.Cells.Item($row, 13).Value = 20
$ws.Cells.Item($row, 14).Value = "Named"

$wb.Save()
